# Newborn outcomes resource file - add new parameter row for
# rr_preterm_death_steroids (antenatal corticosteroids) above the
# "prob_care_seeking_for_complication" row, shifting subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_values")
$ws.Activate()

# Insert a new blank row at row 40 - everything currently at/after row 40
# (including the trailing blank rows) shifts down by one.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new parameter.
$ws.Range("A40").Value = "rr_preterm_death_steroids"
$ws.Range("B40").Value = 0.69
$ws.Range("D40").Value = "Antenatal corticosteroids for accelerating fetal lung maturation for women at risk of preterm birth"

# D40 uses the same bold-font note style as the other "Source" column
# annotations further up the sheet (e.g. D19:D24) rather than the plain
# style used by the "Dummy" notes below it.
$ws.Range("D40").Font.Bold = $true

# Match the author's final selection/active cell.
$ws.Range("D41").Select()
